$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fall1 -> person_laying_but_ok_1 ; truth [200] -> []
$ws.Range("A2").Value = "person_laying_but_ok_1"
$ws.Range("B2").Value = "[]"

# Row 3: fal6_cropped -> simulation_proximus_2_fall_in_line_with_cam_view
$ws.Range("A3").Value = "simulation_proximus_2_fall_in_line_with_cam_view"
$ws.Range("B3").Value = "[135]"
$ws.Range("C3").Value = "[]"
$ws.Range("D3").Value = "[]"

# Row 4: Stroke_Simulation_1 -> person_laying_but_ok_2
$ws.Range("A4").Value = "person_laying_but_ok_2"
$ws.Range("B4").Value = "[]"

# Row 5: FallingAwayFromCamera -> simulation_chantier_2
$ws.Range("A5").Value = "simulation_chantier_2"
$ws.Range("B5").Value = "[700]"
$ws.Range("C5").Value = "[]"
$ws.Range("D5").Value = "[]"

# Row 6: Stroke_Simulation_2 -> young_man_living_2
$ws.Range("A6").Value = "young_man_living_2"
$ws.Range("B6").Value = "[340]"

# Row 7: FallBehindObject -> person_laying_but_ok_3
$ws.Range("A7").Value = "person_laying_but_ok_3"
$ws.Range("B7").Value = "[]"

# Row 8 (new): person_laying_but_ok_4
$ws.Range("A8").Value = "person_laying_but_ok_4"
$ws.Range("B8").Value = "[]"
$ws.Range("C8").Value = "[]"
$ws.Range("D8").Value = "[]"
$ws.Range("E8").Value = "[]"
$ws.Range("F8").Value = "[]"

# Row 9 (new): young_man_watching_tv
$ws.Range("A9").Value = "young_man_watching_tv"
$ws.Range("B9").Value = "[]"
$ws.Range("C9").Value = "[339]"
$ws.Range("D9").Value = "[]"
$ws.Range("E9").Value = "[]"
$ws.Range("F9").Value = "[339]"
